$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 00:42"
$ws.Range("A21").Value = "Brasil"
$ws.Range("B21").Value = 2554
$ws.Range("C21").Value = 307
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 2493
$ws.Range("F21").Value = 18
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = 59
$ws.Range("A22").Value = "Suecia"
$ws.Range("B22").Value = 2526
$ws.Range("C22").Value = 227
$ws.Range("D22").Value = 16
$ws.Range("E22").Value = 2448
$ws.Range("F22").Value = 158
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = 62
$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 558
$ws.Range("C48").Value = 115
$ws.Range("D48").Value = 2
$ws.Range("E48").Value = 548
$ws.Range("F48").Value = 20
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 8
$ws.Range("A49").Value = "Catar"
$ws.Range("B49").Value = 537
$ws.Range("C49").Value = 11
$ws.Range("D49").Value = 41
$ws.Range("E49").Value = 496
$ws.Range("F49").Value = 6
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("A50").Value = "Eslovenia"
$ws.Range("B50").Value = 528
$ws.Range("C50").Value = 48
$ws.Range("D50").Value = 10
$ws.Range("E50").Value = 513
$ws.Range("F50").Value = 14
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 5
$ws.Range("A51").Value = "Peru"
$ws.Range("B51").Value = 480
$ws.Range("C51").Value = 64
$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 470
$ws.Range("F51").Value = 9
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 9
$ws.Range("A52").Value = "Colombia"
$ws.Range("B52").Value = 470
$ws.Range("C52").Value = 92
$ws.Range("D52").Value = 8
$ws.Range("E52").Value = 458
$ws.Range("H52").Value = 4
$ws.Range("A53").Value = "Egipto"
$ws.Range("B53").Value = 456
$ws.Range("C53").Value = 54
$ws.Range("D53").Value = 95
$ws.Range("E53").Value = 340
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 21
$ws.Range("A73").Value = "Uruguay"
$ws.Range("B73").Value = 217
$ws.Range("C73").Value = 28
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 217
$ws.Range("F73").Value = 3
$ws.Range("A74").Value = "Eslovaquia"
$ws.Range("B74").Value = 216
$ws.Range("C74").Value = 12
$ws.Range("D74").Value = 7
$ws.Range("E74").Value = 209
$ws.Range("F74").Value = 2
$ws.Range("H74").Value = 0
$ws.Range("A75").Value = "San Marino"
$ws.Range("B75").Value = 208
$ws.Range("C75").Value = 21
$ws.Range("D75").Value = 4
$ws.Range("F75").Value = 12
$ws.Range("H75").Value = 21
$ws.Range("A76").Value = "Nueva Zelanda"
$ws.Range("B76").Value = 205
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 22
$ws.Range("E76").Value = 183
$ws.Range("F76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("A77").Value = "Costa Rica"
$ws.Range("B77").Value = 201
$ws.Range("C77").Value = 24
$ws.Range("D77").Value = 2
$ws.Range("E77").Value = 197
$ws.Range("F77").Value = 4
$ws.Range("H77").Value = 2
$ws.Range("A78").Value = "Kuwait"
$ws.Range("B78").Value = 195
$ws.Range("C78").Value = 4
$ws.Range("D78").Value = 43
$ws.Range("E78").Value = 152
$ws.Range("F78").Value = 6
$ws.Range("E111").Value = 59
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 1
$ws.Range("A116").Value = "Mauricio"
$ws.Range("C116").Value = 6
$ws.Range("F116").Value = 1
$ws.Range("A117").Value = "Consejo Danes para los Refugiados"
$ws.Range("C117").Value = 3
$ws.Range("F117").Value = 0
$ws.Range("A123").Value = "Mayotte"
$ws.Range("C123").Value = 0
$ws.Range("A124").Value = "Honduras"
$ws.Range("C124").Value = 6
$ws.Range("A131").Value = "Jamaica"
$ws.Range("C131").Value = 5
$ws.Range("D131").Value = 2
$ws.Range("E131").Value = 23
$ws.Range("H131").Value = 1
$ws.Range("A132").Value = "Gibraltar"
$ws.Range("B132").Value = 26
$ws.Range("C132").Value = 11
$ws.Range("D132").Value = 5
$ws.Range("E132").Value = 21
$ws.Range("A133").Value = "Polinesia Francesa"
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 0
$ws.Range("E133").Value = 25
$ws.Range("H133").Value = 0
$ws.Range("A155").Value = "Bermudas"
$ws.Range("C155").Value = 1
$ws.Range("A156").Value = "Seychelles"
$ws.Range("A157").Value = "Dominica"
$ws.Range("C157").Value = 0
$ws.Range("A161").Value = "Gabon"
$ws.Range("A162").Value = "Groenlandia"
$ws.Range("C162").Value = 1
$ws.Range("D162").Value = 2
$ws.Range("E162").Value = 4
$ws.Range("H162").Value = 0
$ws.Range("A163").Value = "Curazao"
$ws.Range("C163").Value = 0
$ws.Range("E163").Value = 3
$ws.Range("H163").Value = 1
$ws.Range("A167").Value = "Guyana"
$ws.Range("D167").Value = 0
$ws.Range("H167").Value = 1
$ws.Range("A168").Value = "Bahamas"
$ws.Range("D168").Value = 1
$ws.Range("H168").Value = 0
$ws.Range("A169").Value = "Guinea"
$ws.Range("A170").Value = "Congo"
$ws.Range("A172").Value = "Suazilandia"
